{"js": "// Remove the stray \"Ver no Jupiter...\" / copyright footer paragraphs\n// (and the blank paragraph that precedes them) that trailed the\n// \"LOM3070: Est\u00e1gio Supervisionado (Requisito)\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the index of the \"Ver no Jupiter...\" paragraph; the blank\n// paragraph immediately preceding it (left over after the\n// \"LOM3070...\" requirement line) is removed along with it and the\n// copyright paragraph that follows.\nlet verIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === targetTexts[0]) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex !== -1) {\n  const toDelete = [];\n  // The blank paragraph right before \"Ver no Jupiter...\" (only delete\n  // it if it is indeed empty, to be safe).\n  if (verIndex - 1 >= 0 && paragraphs.items[verIndex - 1].text.trim() === \"\") {\n    toDelete.push(paragraphs.items[verIndex - 1]);\n  }\n  toDelete.push(paragraphs.items[verIndex]);\n  if (verIndex + 1 < paragraphs.items.length &&\n      paragraphs.items[verIndex + 1].text.trim() === targetTexts[1]) {\n    toDelete.push(paragraphs.items[verIndex + 1]);\n  }\n  toDelete.forEach((p) => p.delete());\n  await context.sync();\n}\n", "ps1": "# Remove the stray \"Ver no Jupiter...\" / copyright footer paragraphs\n# (and the blank paragraph that precedes them) that trailed the\n# \"LOM3070: Est\u00e1gio Supervisionado (Requisito)\" paragraph.\n\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [string][char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$verIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text.Trim()\n    if ($t -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -gt 0) {\n    $toDelete = @()\n\n    # Blank paragraph immediately before \"Ver no Jupiter...\" (only if empty).\n    if ($verIndex - 1 -ge 1) {\n        $prevText = $paras.Item($verIndex - 1).Range.Text.Trim()\n        if ($prevText -eq \"\") {\n            $toDelete += ($verIndex - 1)\n        }\n    }\n\n    $toDelete += $verIndex\n\n    if ($verIndex + 1 -le $count) {\n        $nextText = $paras.Item($verIndex + 1).Range.Text.Trim()\n        if ($nextText -eq $copyrightText) {\n            $toDelete += ($verIndex + 1)\n        }\n    }\n\n    # Delete from highest index to lowest so earlier indices stay valid.\n    $sorted = $toDelete | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
